# Insert a new price-record row for Perejil (Feria Lagunitas de Puerto Montt)
# above the current row 123. Excel shifts rows 123..201 down to 124..202,
# which grows the sheet's used range from A1:R201 to A1:R202.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(123).Insert()

# Populate the newly-inserted row 123 with the new weekly record.
$ws.Range("A123").Value = 4
$ws.Range("B123").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C123").Value = "Los Lagos"
$ws.Range("D123").Value = 44582
$ws.Range("E123").Value = 10
$ws.Range("F123").Value = 100112044
$ws.Range("G123").Value = "Perejil"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 180
$ws.Range("K123").Value = 5000
$ws.Range("L123").Value = 5000
$ws.Range("M123").Value = 5000
$ws.Range("N123").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O123").Value = "Región Metropolitana"
$ws.Range("P123").Value = 1667
$ws.Range("Q123").Value = 3
$ws.Range("R123").Value = "Hortaliza"
